$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.450.98'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.161.09'
$ws.Range("E3").Value = '  +3.18%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.62'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.30'
$ws.Range("E7").Value = '  +4.38%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0861'
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.05'
$ws.Range("E12").Value = '  +4.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.480.60'
$ws.Range("E13").Value = '  +3.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.27'
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.813'
$ws.Range("E15").Value = '  +0.62%  '
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.173.72'
$ws.Range("E17").Value = '  +3.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.407.59'
$ws.Range("E18").Value = '  +1.63%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.80'
$ws.Range("E19").Value = '  -0.16%  '
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0854'
$ws.Range("E21").Value = '  +1.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.69'
$ws.Range("E22").Value = '  +1.51%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.52'
$ws.Range("E24").Value = '  +5.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  +0.69%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.67'
$ws.Range("E26").Value = '  +1.28%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.38'
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.89'
$ws.Range("E29").Value = '  +2.90%  '
$ws.Range("E31").Value = '  +7.10%  '
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.63'
$ws.Range("E33").Value = '  +2.27%  '
$ws.Range("E34").Value = '  -0.32%  '
$ws.Range("E35").Value = '  +8.88%  '
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '104.21'
$ws.Range("E40").Value = '  +2.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0231'
$ws.Range("E41").Value = '  +0.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.83'
$ws.Range("E42").Value = '  -1.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.539.93'
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.91'
$ws.Range("E45").Value = '  +3.05%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0925'
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("B47").Value = 'HuobiToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.82'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("E48").Value = '  +5.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.21'
$ws.Range("E49").Value = '  +2.36%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.364.12'
$ws.Range("E50").Value = '  +3.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.96'
$ws.Range("E51").Value = '  +0.10%  '
